$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.827.37"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "'3.763.11"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'621.34"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'181.09"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").Value = "'3.761.51"
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("E10").Value = "  +4.21%  "
$ws.Range("D11").Value = "'6.31"
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "'41.45"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "'4.379.20"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("D16").Value = "'3.759.56"
$ws.Range("E16").Value = "  +2.79%  "
$ws.Range("D17").Value = "'69.892.60"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'7.64"
$ws.Range("E19").Value = "  +1.77%  "
$ws.Range("D20").Value = "'16.79"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("D21").Value = "'508.42"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("D23").Value = "'0.729"
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("D24").Value = "'2.52"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("D25").Value = "'87.32"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "'13.19"
$ws.Range("E26").Value = "  -1.75%  "
$ws.Range("D27").Value = "'11.14"
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("D28").Value = "'0.0000137"
$ws.Range("E28").Value = "  +25.18%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'2.52"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").Value = "'2.90"
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("D32").Value = "'7.90"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("D34").Value = "'0.116"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("D37").Value = "'6.24"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("D38").Value = "'0.338"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "'0.132"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("D40").Value = "'2.12"
$ws.Range("E40").Value = "  -2.77%  "
$ws.Range("D41").Value = "'50.26"
$ws.Range("D42").Value = "'45.74"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "'428.51"
$ws.Range("E43").Value = "  +2.21%  "
$ws.Range("D44").Value = "'8.75"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Value = "'2.85"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("D46").Value = "'3.007.64"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("D48").Value = "'27.48"
$ws.Range("E48").Value = "  -4.84%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'136.65"
$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  +1.41%  "
